$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M38").Value = -572.1428999999999
$ws.Range("K38").Value = 944.1428999999999
$ws.Range("J38").Value = 6000
$ws.Range("I38").Value = 314.7143
$ws.Range("L38").Value = 18000
$ws.Range("H38").Value = 1025.375
$ws.Range("M97").ClearContents()
$ws.Range("H97").Value = 4495.3335
$ws.Range("I97").Value = 0
$ws.Range("L97").Value = 13486.0005
$ws.Range("N97").Value = -14478.0005
$ws.Range("K97").Value = 0
$ws.Range("J97").Value = 4495.3335
$ws.Range("J100").Value = 2250
$ws.Range("I100").Value = 126995.125
$ws.Range("H100").Value = 113134.555
$ws.Range("M100").Value = -126454.125
$ws.Range("K100").Value = 126995.125
$ws.Range("L100").Value = 2250
$ws.Range("N100").Value = -3332
$ws.Range("H137").Value = 13090.35
$ws.Range("K137").Value = 29078.316
$ws.Range("I137").Value = 9692.772000000001
$ws.Range("M137").Value = -26528.316
$ws.Range("J137").Value = 14913.439
$ws.Range("L137").Value = 44740.317

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M32").Value = -2858.9312
$ws.Range("H32").Value = 4493.4443
$ws.Range("N32").Value = -10650
$ws.Range("I32").Value = 3145.9312
$ws.Range("L32").Value = 10076
$ws.Range("K32").Value = 3145.9312
$ws.Range("J32").Value = 10076
$ws.Range("K61").Value = 3131
$ws.Range("I61").Value = 3131
$ws.Range("J61").Value = 16064.667
$ws.Range("L61").Value = 16064.667
$ws.Range("H61").Value = 6364.4165
$ws.Range("M61").Value = -2919
$ws.Range("N61").Value = -16488.667
$ws.Range("K80").Value = 53000
$ws.Range("N80").Value = -45401
$ws.Range("H80").Value = 44775.715
$ws.Range("L80").Value = 43405
$ws.Range("M80").Value = -52002
$ws.Range("J80").Value = 43405
$ws.Range("I80").Value = 53000
$ws.Range("M83").Value = -154008
$ws.Range("I83").Value = 53000
$ws.Range("H83").Value = 44775.715
$ws.Range("J83").Value = 43405
$ws.Range("L83").Value = 130215
$ws.Range("N83").Value = -140199
$ws.Range("K83").Value = 159000
$ws.Range("M102").Value = -607.6667000000002
$ws.Range("K102").Value = 2229.6667
$ws.Range("I102").Value = 2229.6667
$ws.Range("H102").Value = 2229.6667
$ws.Range("L102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("M132").Value = -41319.764
$ws.Range("J132").Value = 29669.428
$ws.Range("H132").Value = 19007
$ws.Range("K132").Value = 43849.764
$ws.Range("I132").Value = 14616.588
$ws.Range("L132").Value = 89008.284
$ws.Range("K136").Value = 9393
$ws.Range("H136").Value = 6364.4165
$ws.Range("I136").Value = 3131
$ws.Range("N136").Value = -53294.001
$ws.Range("J136").Value = 16064.667
$ws.Range("L136").Value = 48194.001
$ws.Range("M136").Value = -6843

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N42").Value = -320655
$ws.Range("L42").Value = 319999
$ws.Range("K42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 319999
$ws.Range("H42").Value = 319999
$ws.Range("N86").Value = -4563.5386
$ws.Range("I86").Value = 478300.25
$ws.Range("K86").Value = 478300.25
$ws.Range("H86").Value = 296306.84
$ws.Range("L86").Value = 2317.5386
$ws.Range("M86").Value = -477177.25
$ws.Range("J86").Value = 2317.5386
$ws.Range("M89").Value = -2385885.25
$ws.Range("N89").Value = -22819.693
$ws.Range("K89").Value = 2391501.25
$ws.Range("I89").Value = 478300.25
$ws.Range("L89").Value = 11587.693
$ws.Range("J89").Value = 2317.5386
$ws.Range("H89").Value = 296306.84
$ws.Range("L94").Value = 999.5
$ws.Range("I94").Value = 835.7
$ws.Range("J94").Value = 999.5
$ws.Range("M94").Value = -384.7
$ws.Range("N94").Value = -1901.5
$ws.Range("K94").Value = 835.7
$ws.Range("H94").Value = 863
$ws.Range("J99").Value = 3375
$ws.Range("L99").Value = 3375
$ws.Range("N99").Value = -6371
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("I99").Value = 0
$ws.Range("H99").Value = 3375
$ws.Range("H105").Value = 5605.472
$ws.Range("I105").Value = 7341.913
$ws.Range("M105").Value = -5594.913
$ws.Range("N105").Value = -6027.3076
$ws.Range("J105").Value = 2533.3076
$ws.Range("L105").Value = 2533.3076
$ws.Range("K105").Value = 7341.913

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K16").Value = 2982.6155
$ws.Range("J16").Value = 2409
$ws.Range("I16").Value = 2982.6155
$ws.Range("M16").Value = -2695.6155
$ws.Range("H16").Value = 2823.2778
$ws.Range("L16").Value = 2409
$ws.Range("L31").Value = 1894.6666
$ws.Range("H31").Value = 1753.3062
$ws.Range("K31").Value = 1461.75
$ws.Range("J31").Value = 1894.6666
$ws.Range("I31").Value = 1461.75
$ws.Range("M31").Value = -1166.75
$ws.Range("J34").Value = 1894.6666
$ws.Range("K34").Value = 1461.75
$ws.Range("H34").Value = 1753.3062
$ws.Range("M34").Value = -1259.75
$ws.Range("I34").Value = 1461.75
$ws.Range("L34").Value = 1894.6666
$ws.Range("L113").Value = 2409
$ws.Range("J113").Value = 2409
$ws.Range("M113").Value = -812.6154999999999
$ws.Range("H113").Value = 2823.2778
$ws.Range("I113").Value = 2982.6155
$ws.Range("K113").Value = 2982.6155

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J5").Value = 5221.154
$ws.Range("I5").Value = 2981.1667
$ws.Range("K5").Value = 8943.500100000001
$ws.Range("L5").Value = 15663.462
$ws.Range("H5").Value = 4145.96
$ws.Range("M5").Value = -8831.500100000001
$ws.Range("M32").ClearContents()
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("J60").Value = 2165.6667
$ws.Range("N60").Value = -6999.000100000001
$ws.Range("L60").Value = 6497.000100000001
$ws.Range("H60").Value = 2165.6667
$ws.Range("K60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N86").ClearContents()
$ws.Range("I86").Value = 899.5
$ws.Range("K86").Value = 2698.5
$ws.Range("H86").Value = 899.5
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1512.5
$ws.Range("J86").Value = 0
$ws.Range("M89").Value = -2167.5
$ws.Range("N89").ClearContents()
$ws.Range("K89").Value = 8095.5
$ws.Range("I89").Value = 899.5
$ws.Range("L89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("H89").Value = 899.5
$ws.Range("H105").Value = 12500
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 12500
$ws.Range("L105").Value = 37500
$ws.Range("K105").Value = 0
$ws.Range("N105").Value = -42742
$ws.Range("H131").Value = 3125.0527
$ws.Range("K131").Value = 3927
$ws.Range("L131").Value = 12553.2495
$ws.Range("M131").Value = 1113
$ws.Range("J131").Value = 4184.4165
$ws.Range("N131").Value = -22633.2495
$ws.Range("I131").Value = 1309
$ws.Range("I135").Value = 2981.1667
$ws.Range("K135").Value = 26830.5003
$ws.Range("L135").Value = 46990.38600000001
$ws.Range("H135").Value = 4145.96
$ws.Range("J135").Value = 5221.154
$ws.Range("M135").Value = -24295.5003
$ws.Range("I139").Value = 1739.3334
$ws.Range("M139").Value = -78.0002000000004
$ws.Range("K139").Value = 5218.0002
$ws.Range("L139").Value = 11997
$ws.Range("N139").Value = -22277
$ws.Range("H139").Value = 2869.1667
$ws.Range("J139").Value = 3999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I2").Value = 87.71429000000001
$ws.Range("N2").Value = -434.5
$ws.Range("M2").Value = 25.28570999999999
$ws.Range("H2").Value = 143.46153
$ws.Range("L2").Value = 208.5
$ws.Range("K2").Value = 87.71429000000001
$ws.Range("J2").Value = 208.5
$ws.Range("L127").Value = 40220.5
$ws.Range("K127").Value = 0
$ws.Range("J127").Value = 40220.5
$ws.Range("I127").Value = 0
$ws.Range("N127").Value = -50140.5
$ws.Range("H127").Value = 40220.5
$ws.Range("M132").Value = -19.0001000000002
$ws.Range("J132").Value = 1200
$ws.Range("H132").Value = 937.25
$ws.Range("K132").Value = 2549.0001
$ws.Range("I132").Value = 849.6667
$ws.Range("L132").Value = 3600

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M7").Value = -1130.1818
$ws.Range("H7").Value = 1478.875
$ws.Range("K7").Value = 1242.1818
$ws.Range("J7").Value = 1999.6
$ws.Range("I7").Value = 1242.1818
$ws.Range("N7").Value = -2223.6
$ws.Range("L7").Value = 1999.6
$ws.Range("K16").Value = 1579.7142
$ws.Range("J16").Value = 3248
$ws.Range("I16").Value = 1579.7142
$ws.Range("N16").Value = -3588
$ws.Range("M16").Value = -1409.7142
$ws.Range("H16").Value = 2186.3635
$ws.Range("L16").Value = 3248
$ws.Range("I22").Value = 1595.3334
$ws.Range("L22").Value = 4718.75
$ws.Range("J22").Value = 4718.75
$ws.Range("K22").Value = 1595.3334
$ws.Range("H22").Value = 2556.3845
$ws.Range("M22").Value = -1300.3334
$ws.Range("I27").Value = 1595.3334
$ws.Range("H27").Value = 2556.3845
$ws.Range("J27").Value = 4718.75
$ws.Range("M27").Value = -1488.3334
$ws.Range("L27").Value = 4718.75
$ws.Range("K27").Value = 1595.3334
$ws.Range("H46").Value = 1651.2354
$ws.Range("M46").Value = -864.7778000000001
$ws.Range("K46").Value = 1052.7778
$ws.Range("I46").Value = 1052.7778
$ws.Range("L46").Value = 2324.5
$ws.Range("J46").Value = 2324.5
$ws.Range("J47").Value = 0
$ws.Range("I47").Value = 59
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H47").Value = 59
$ws.Range("M47").Value = 431
$ws.Range("K47").Value = 59
$ws.Range("L52").Value = 0
$ws.Range("H52").Value = 59
$ws.Range("I52").Value = 59
$ws.Range("M52").Value = 174
$ws.Range("K52").Value = 59
$ws.Range("N52").ClearContents()
$ws.Range("J52").Value = 0
$ws.Range("H55").Value = 649.13336
$ws.Range("N55").Value = -559
$ws.Range("M55").Value = -694.2
$ws.Range("J55").Value = 213
$ws.Range("I55").Value = 867.2
$ws.Range("L55").Value = 213
$ws.Range("K55").Value = 867.2
$ws.Range("J100").Value = 1236900
$ws.Range("I100").Value = 7665.6665
$ws.Range("H100").Value = 929591.4399999999
$ws.Range("M100").Value = -7124.6665
$ws.Range("K100").Value = 7665.6665
$ws.Range("L100").Value = 1236900
$ws.Range("N100").Value = -1237982
$ws.Range("N126").Value = -10938.8
$ws.Range("H126").Value = 1478.875
$ws.Range("I126").Value = 1242.1818
$ws.Range("L126").Value = 5998.799999999999
$ws.Range("M126").Value = -1256.5454
$ws.Range("J126").Value = 1999.6
$ws.Range("K126").Value = 3726.5454
$ws.Range("M132").Value = -8791526.600000001
$ws.Range("J132").Value = 8459.5
$ws.Range("H132").Value = 2532776
$ws.Range("K132").Value = 8794056.600000001
$ws.Range("I132").Value = 2931352.2
$ws.Range("L132").Value = 25378.5
$ws.Range("K136").Value = 25648857
$ws.Range("H136").Value = 5054730
$ws.Range("I136").Value = 8549619
$ws.Range("N136").Value = -24770.6661
$ws.Range("J136").Value = 6556.8887
$ws.Range("L136").Value = 19670.6661
$ws.Range("M136").Value = -25646307

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J100").Value = 749.6667
$ws.Range("I100").Value = 749.1
$ws.Range("H100").Value = 749.2308
$ws.Range("M100").Value = -957.2
$ws.Range("K100").Value = 1498.2
$ws.Range("L100").Value = 1499.3334
$ws.Range("N100").Value = -2581.3334
